$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was refreshed by one day
# (from 45181 = 2023-09-12 to 45182 = 2023-09-13) for every data row (2..200).
for ($r = 2; $r -le 200; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
